# Update the "Expertise" column (G) on the "6 Elmer Seawood" sheet with
# revised L/M/H/blank ratings, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("6 Elmer Seawood")

# Row -> new Expertise value ("" clears the cell)
$updates = @{
    2  = ""
    3  = ""
    5  = "L"
    6  = "L"
    8  = "L"
    9  = "M"
    11 = ""
    13 = ""
    14 = ""
    15 = ""
    18 = "L"
    20 = "L"
    21 = "M"
    22 = ""
    23 = ""
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    if ($value -eq "") {
        $ws.Range("G$row").ClearContents()
    } else {
        $ws.Range("G$row").Value = $value
    }
}

# Move the current selection on the sheet
[void]$ws.Activate()
[void]$ws.Range("A10:K10").Select()
